$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl21b"
$ws.Range("C2").Value = "Ccr7"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.004739
$ws.Range("H2").Value = 0.014217
$ws.Range("I2").Value = 0.02588570741885795
$ws.Range("J2").Value = 0.02588570741885795
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.186073
$ws.Range("N2").Value = 0.558219
$ws.Range("O2").Value = 0.08082280080999586
$ws.Range("P2").Value = 0.08082280080999586
$ws.Range("Q2").Value = 0.0008817999470000001
$ws.Range("R2").Value = 0.007936199523000001
$ws.Range("S2").Value = 0.002092155374540188
$ws.Range("T2").Value = 0.002092155374540188

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl21b"
$ws.Range("C3").Value = "Ccr7"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.004739
$ws.Range("H3").Value = 0.014217
$ws.Range("I3").Value = 0.02588570741885795
$ws.Range("J3").Value = 0.02588570741885795
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.116161
$ws.Range("N3").Value = 6.348483
$ws.Range("O3").Value = 0.9191771991900042
$ws.Range("P3").Value = 0.9191771991900041
$ws.Range("Q3").Value = 0.010028486979
$ws.Range("R3").Value = 0.090256382811
$ws.Range("S3").Value = 0.02379355204431776
$ws.Range("T3").Value = 0.02379355204431776

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Ccl21b"
$ws.Range("C4").Value = "Ccr7"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.134289
$ws.Range("H4").Value = 0.402867
$ws.Range("I4").Value = 0.733523056250478
$ws.Range("J4").Value = 0.733523056250478
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.186073
$ws.Range("N4").Value = 0.558219
$ws.Range("O4").Value = 0.08082280080999586
$ws.Range("P4").Value = 0.08082280080999586
$ws.Range("Q4").Value = 0.02498755709700001
$ws.Range("R4").Value = 0.224888013873
$ws.Range("S4").Value = 0.05928538786487177
$ws.Range("T4").Value = 0.05928538786487177

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl21b"
$ws.Range("C5").Value = "Ccr7"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.134289
$ws.Range("H5").Value = 0.402867
$ws.Range("I5").Value = 0.733523056250478
$ws.Range("J5").Value = 0.733523056250478
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.116161
$ws.Range("N5").Value = 6.348483
$ws.Range("O5").Value = 0.9191771991900042
$ws.Range("P5").Value = 0.9191771991900041
$ws.Range("Q5").Value = 0.2841771445290001
$ws.Range("R5").Value = 2.557594300761
$ws.Range("S5").Value = 0.6742376683856063
$ws.Range("T5").Value = 0.6742376683856062

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Ccl21b"
$ws.Range("C6").Value = "Ccr7"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.04404599999999999
$ws.Range("H6").Value = 0.132138
$ws.Range("I6").Value = 0.2405912363306641
$ws.Range("J6").Value = 0.2405912363306641
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.186073
$ws.Range("N6").Value = 0.558219
$ws.Range("O6").Value = 0.08082280080999586
$ws.Range("P6").Value = 0.08082280080999586
$ws.Range("Q6").Value = 0.008195771358
$ws.Range("R6").Value = 0.07376194222199998
$ws.Range("S6").Value = 0.0194452575705839
$ws.Range("T6").Value = 0.0194452575705839

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Ccl21b"
$ws.Range("C7").Value = "Ccr7"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.04404599999999999
$ws.Range("H7").Value = 0.132138
$ws.Range("I7").Value = 0.2405912363306641
$ws.Range("J7").Value = 0.2405912363306641
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.116161
$ws.Range("N7").Value = 6.348483
$ws.Range("O7").Value = 0.9191771991900042
$ws.Range("P7").Value = 0.9191771991900041
$ws.Range("Q7").Value = 0.09320842740599999
$ws.Range("R7").Value = 0.8388758466539998
$ws.Range("S7").Value = 0.2211459787600802
$ws.Range("T7").Value = 0.2211459787600802

